# Updates the crypto price/volume table to the latest scrape.
# All values are written as TEXT (leading "'" forces Excel to store them
# as strings rather than auto-coercing numeric-looking text like "0.999"
# or "8.80" into numbers), matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.521.10'
$ws.Range('E2').Value = '''  -6.05%  '
$ws.Range('D3').Value = '''3.526.32'
$ws.Range('E3').Value = '''  -1.73%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '''  -0.06%  '
$ws.Range('D5').Value = '''391.37'
$ws.Range('E5').Value = '''  -6.22%  '
$ws.Range('D6').Value = '''121.25'
$ws.Range('E6').Value = '''  -7.11%  '
$ws.Range('D7').Value = '''3.514.31'
$ws.Range('E7').Value = '''  -1.74%  '
$ws.Range('D8').Value = '''0.588'
$ws.Range('E8').Value = '''  -9.58%  '
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '''  +0.03%  '
$ws.Range('D10').Value = '''0.679'
$ws.Range('E10').Value = '''  -11.90%  '
$ws.Range('D11').Value = '''0.152'
$ws.Range('E11').Value = '''  -14.43%  '
$ws.Range('D12').Value = '''0.0000339'
$ws.Range('E12').Value = '''  -0.94%  '
$ws.Range('D13').Value = '''38.76'
$ws.Range('E13').Value = '''  -8.60%  '
$ws.Range('D14').Value = '''4.073.09'
$ws.Range('E14').Value = '''  -1.56%  '
$ws.Range('D15').Value = '''9.24'
$ws.Range('E15').Value = '''  -7.72%  '
$ws.Range('E16').Value = '''  -3.23%  '
$ws.Range('D17').Value = '''3.523.14'
$ws.Range('E17').Value = '''  -1.90%  '
# Row 18/19 swapped: Uniswap now ranks above Chainlink
$ws.Range('B18').Value = '''Uniswap'
$ws.Range('C18').Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '''12.72'
$ws.Range('E18').Value = '''  +2.59%  '
$ws.Range('B19').Value = '''Chainlink'
$ws.Range('C19').Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''18.78'
$ws.Range('E19').Value = '''  -8.24%  '
$ws.Range('D20').Value = '''63.355.53'
$ws.Range('E20').Value = '''  -5.94%  '
$ws.Range('D21').Value = '''1.02'
$ws.Range('E21').Value = '''  -11.43%  '
$ws.Range('D22').Value = '''394.15'
$ws.Range('E22').Value = '''  -14.10%  '
$ws.Range('D23').Value = '''13.92'
$ws.Range('E23').Value = '''  +3.18%  '
$ws.Range('D24').Value = '''80.98'
$ws.Range('E24').Value = '''  -8.41%  '
$ws.Range('D25').Value = '''2.85'
$ws.Range('E25').Value = '''  -7.97%  '
$ws.Range('D26').Value = '''33.24'
$ws.Range('E26').Value = '''  -5.07%  '
$ws.Range('D27').Value = '''5.15'
$ws.Range('E27').Value = '''  +5.90%  '
$ws.Range('E28').Value = '''  -11.80%  '
$ws.Range('D29').Value = '''8.80'
$ws.Range('E29').Value = '''  -13.11%  '
$ws.Range('D30').Value = '''11.94'
$ws.Range('E30').Value = '''  -3.81%  '
$ws.Range('D31').Value = '''2.61'
$ws.Range('E31').Value = '''  -6.75%  '
$ws.Range('E32').Value = '''  -6.15%  '
$ws.Range('D33').Value = '''6.85'
$ws.Range('E33').Value = '''  -8.26%  '
$ws.Range('D34').Value = '''0.151'
$ws.Range('E34').Value = '''  -7.28%  '
$ws.Range('E35').Value = '''  +0.09%  '
$ws.Range('E36').Value = '''  -12.31%  '
$ws.Range('D37').Value = '''54.07'
$ws.Range('E37').Value = '''  -4.48%  '
$ws.Range('D38').Value = '''0.0438'
$ws.Range('E38').Value = '''  -11.56%  '
$ws.Range('D39').Value = '''0.996'
$ws.Range('E39').Value = '''  -0.14%  '
$ws.Range('E40').Value = '''  +16.33%  '
$ws.Range('D41').Value = '''0.0₃0632'
$ws.Range('E41').Value = '''  -12.50%  '
$ws.Range('E42').Value = '''  -10.65%  '
$ws.Range('E43').Value = '''  +14.23%  '
$ws.Range('D44').Value = '''141.39'
$ws.Range('E44').Value = '''  -5.14%  '
# Rows 45-47 reshuffled: LidoDAOToken, WEMIXToken, Stacks
$ws.Range('B45').Value = '''LidoDAOToken'
$ws.Range('C45').Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D45').Value = '''3.08'
$ws.Range('E45').Value = '''  -5.98%  '
$ws.Range('B46').Value = '''WEMIXToken'
$ws.Range('C46').Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = '''2.49'
$ws.Range('E46').Value = '''  -8.98%  '
$ws.Range('B47').Value = '''Stacks'
$ws.Range('C47').Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '''2.70'
$ws.Range('E47').Value = '''  -11.21%  '
$ws.Range('D48').Value = '''1.95'
$ws.Range('E48').Value = '''  -1.66%  '
$ws.Range('D49').Value = '''24.64'
$ws.Range('E49').Value = '''  +13.27%  '
$ws.Range('D50').Value = '''4.03'
$ws.Range('E50').Value = '''  -6.82%  '
$ws.Range('E51').Value = '''  -10.86%  '
